$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FE): clear old B/D estimates, update C to the revised value
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0.2
$ws.Range("D2").ClearContents()

# Row 3 (FE+Disg): clear old B/D estimates, update C to the revised value
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 0.2
$ws.Range("D3").ClearContents()

# Row 4 (FE+Disg+Var): clear old B/D estimates, update C to the revised value
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = 0.2
$ws.Range("D4").ClearContents()

# Match the saved workbook's selection / cursor position
$null = $ws.Range("G16").Select()

# Match the saved workbook's page margins (values are in points; 72pt = 1in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
